$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2/R2 to the nearest integer
$ws.Range("Q2").Value = 602196
$ws.Range("R2").Value = 6555866

# Clear the start time (Z2) and end time (AB2) cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
